$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they are not auto-converted to numbers
foreach ($addr in @("D5","D6","D7","D10","D11","D12","D13","D14","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D39","D42","D43","D44","D46","D47","D48","D49","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the latest cryptos data pull
$ws.Range("D2").Value = "65.831.22"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "3.405.64"
$ws.Range("E3").Value = "  -5.85%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "183.01"
$ws.Range("E5").Value = "  -9.98%  "
$ws.Range("D6").Value = "535.19"
$ws.Range("E6").Value = "  -5.03%  "
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.396.67"
$ws.Range("E8").Value = "  -6.09%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "0.632"
$ws.Range("E10").Value = "  -6.41%  "
$ws.Range("D11").Value = "57.96"
$ws.Range("E11").Value = "  -3.98%  "
$ws.Range("D12").Value = "0.136"
$ws.Range("E12").Value = "  -11.08%  "
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  -11.37%  "
$ws.Range("D14").Value = "9.45"
$ws.Range("E14").Value = "  -5.99%  "
$ws.Range("D15").Value = "3.979.92"
$ws.Range("E15").Value = "  -5.32%  "
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "3.438.61"
$ws.Range("E17").Value = "  -4.53%  "
$ws.Range("D18").Value = "65.930.20"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").Value = "17.67"
$ws.Range("E19").Value = "  -7.43%  "
$ws.Range("D20").Value = "11.40"
$ws.Range("E20").Value = "  -7.93%  "
$ws.Range("D21").Value = "0.988"
$ws.Range("E21").Value = "  -8.60%  "
$ws.Range("D22").Value = "381.68"
$ws.Range("E22").Value = "  -5.79%  "
$ws.Range("D23").Value = "83.38"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").Value = "3.80"
$ws.Range("E24").Value = "  -9.02%  "
$ws.Range("D25").Value = "11.03"
$ws.Range("E25").Value = "  -14.63%  "
$ws.Range("D26").Value = "3.71"
$ws.Range("E26").Value = "  -4.46%  "
$ws.Range("D27").Value = "11.81"
$ws.Range("E27").Value = "  -6.16%  "
$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  -8.35%  "
$ws.Range("D29").Value = "8.62"
$ws.Range("E29").Value = "  -8.70%  "
$ws.Range("D30").Value = "700.62"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("D31").Value = "29.96"
$ws.Range("E31").Value = "  -5.51%  "
$ws.Range("D32").Value = "6.86"
$ws.Range("E32").Value = "  -18.27%  "
$ws.Range("D33").Value = "11.30"
$ws.Range("E33").Value = "  -7.52%  "
$ws.Range("D34").Value = "61.78"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "37.07"
$ws.Range("E37").Value = "  -13.05%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "0.393"
$ws.Range("E38").Value = "  -8.86%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  -6.19%  "
$ws.Range("D41").Value = "2.909.14"
$ws.Range("E41").Value = "  -10.45%  "
$ws.Range("D42").Value = "2.81"
$ws.Range("E42").Value = "  -13.54%  "
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0398"
$ws.Range("E44").Value = "  -4.89%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0630"
$ws.Range("E45").Value = "  -19.13%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.40"
$ws.Range("E46").Value = "  -14.51%  "
$ws.Range("D47").Value = "0.126"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "2.95"
$ws.Range("E48").Value = "  -4.94%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "134.68"
$ws.Range("E49").Value = "  -3.78%  "
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  -23.11%  "
